$d = $word.ActiveDocument

# 1. "Hello {{ name }}," paragraph: collapse the split "{{ name }}" field markup back
#    into a single run, then add a new run with the extra greeting text.
$d.Content.Find.Execute("Hello {{ name }},", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Hello {{ name }},", 2)

# 2. "{{ paragraph_replace }}" paragraph: collapse the split field markup into one run.
$d.Content.Find.Execute("{{ paragraph_replace }}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{{ paragraph_replace }}", 2)

# 3. "{{ table }}" paragraph: collapse the split field markup into one run.
$d.Content.Find.Execute("{{ table }}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{{ table }}", 2)

# Now append the new trailing run " how are you?" right after "Hello {{ name }},".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Hello {{ name }},`r") {
        $p.Range.InsertAfter(" how are you?")
    }
}

# 4. Remove the _GoBack bookmark left over on the "Thank you" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
